$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 469.0027611255646
$ws.Range("B2").Value = 7.816712685426077
$ws.Range("C2").Value = 66.99849888256618
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 0.895517115916413
$ws.Range("F2").Value = "0.9 epochs/min"
